$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Next period (release date)" for Annual Population Survey rows (D2 and D4)
# from "Apr 2022 - Mar 2023 (11/06/23)" to "Apr 2022 - Mar 2023 (11/07/23)"
$ws.Range("D2").Value = "Apr 2022 - Mar 2023 (11/07/23)"
$ws.Range("D4").Value = "Apr 2022 - Mar 2023 (11/07/23)"

# Match the resulting cell selection left behind after the edit
$ws.Activate()
$ws.Range("D5").Select()
